$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update date, initials, call off, room number; clear comments cell
$ws.Range("A2").Value = 43068
$ws.Range("B2").Value = "MS"
$ws.Range("C2").Value = "NO"
$ws.Range("D2").Value = "2"
$ws.Range("E2").ClearContents()

# Rows 3-12: clear all debugging/sample data contents (formatting/styles stay intact)
$ws.Range("A3:F12").ClearContents()
